# Update the "want to go" counts (column F) for three events that appear
# on both the "展览" sheet and the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 323
$wsExhibit.Range("F4").Value = 2855
$wsExhibit.Range("F6").Value = 604

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 323
$wsAll.Range("F6").Value = 2855
$wsAll.Range("F8").Value = 604
